# Groups2.xlsx: a 7th child joined the carpool. The per-child table (rows
# 6-13 after this edit) gets re-sorted by remaining distance/time (col H,
# descending) and renumbered; the "school"/"cost"/"time" summary rows shift
# down by one to make room for the new child row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every first/last name in this sheet carries a trailing " " + NBSP marker
# (e.g. stored "Kandis" is actually "Kandis" + space + U+00A0) - preserve it.
$nbsp = [char]0x00A0
function NameText($s) {
    return "$s " + $nbsp
}

# Insert a new row right before the old "school" row (row 13), pushing
# school/cost/time down one row (13->14, 14->15, 15->16).
$ws.Rows.Item(13).Insert()

# nChildren: 7 -> 8
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "8"

# Children table rows 6-13
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "0"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "9"
$ws.Range("C6").Value = (NameText "Letha")
$ws.Range("D6").Value = (NameText "Stephenie")
$ws.Range("E6").Value = "-6.44,9.6"
$ws.Range("F6").Value = "Sibyl(mother): 0567328221"
$ws.Range("G6").Value = "7:00:00"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "27.0"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "1"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "13"
$ws.Range("C7").Value = (NameText "Fay")
$ws.Range("D7").Value = (NameText "Emilee")
$ws.Range("E7").Value = "-8.3,6.81"
$ws.Range("F7").Value = "Sheri(mother): 0516797453"
$ws.Range("G7").Value = "7:05:00"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "22.0"

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "7"
$ws.Range("C8").Value = (NameText "Wyatt")
$ws.Range("D8").Value = (NameText "Willette")
$ws.Range("E8").Value = "-7.84,3.24"
$ws.Range("F8").Value = "Antionette(father): 0557331799"
$ws.Range("G8").Value = "7:10:00"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "17.0"

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "3"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "1"
$ws.Range("C9").Value = (NameText "Corene")
$ws.Range("D9").Value = (NameText "Myra")
$ws.Range("E9").Value = "-7.45,3.53"
$ws.Range("F9").Value = "Georgie(mother): 0544823581"
$ws.Range("G9").Value = "7:11:00"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "16.0"

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "4"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "17"
$ws.Range("C10").Value = (NameText "Britta")
$ws.Range("D10").Value = (NameText "Jamel")
$ws.Range("E10").Value = "-5.94,3.44"
$ws.Range("F10").Value = "Albertine(father): 0574981040"
$ws.Range("G10").Value = "7:13:00"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "14.0"

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "5"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "15"
$ws.Range("C11").Value = (NameText "Nubia")
$ws.Range("D11").Value = (NameText "Royce")
$ws.Range("E11").Value = "-3.23,2.78"
$ws.Range("F11").Value = "Augustus(father): 0517389040"
$ws.Range("G11").Value = "7:17:00"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "10.0"

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "6"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "19"
$ws.Range("C12").Value = (NameText "Jeanine")
$ws.Range("D12").Value = (NameText "Janee")
$ws.Range("E12").Value = "-2.97,1.58"
$ws.Range("F12").Value = "Teresa(mother): 0517627420"
$ws.Range("G12").Value = "7:19:00"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "8.0"

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "7"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "10"
$ws.Range("C13").Value = (NameText "Demetra")
$ws.Range("D13").Value = (NameText "Francene")
$ws.Range("E13").Value = "-3.45,-0.28"
$ws.Range("F13").Value = "Dorian(mother): 0534328089"
$ws.Range("G13").Value = "7:22:00"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "5.0"

# "school" summary row, now at row 14
$ws.Range("A14").Value = "school"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "3"
$ws.Range("C14").Value = "Ironiah"
$ws.Range("D14").Value = "mySchool"
$ws.Range("E14").Value = "0,0"
$ws.Range("F14").Value = "Shir(secretary): 0523345098"
$ws.Range("G14").Value = "7:27:00"

# "cost" row, now at row 15 (value unchanged)
$ws.Range("A15").Value = "cost"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "25"

# "time" row, now at row 16
$ws.Range("A16").Value = "time"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "27.0"
